$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Confirmed Cases")

# New confirmed-case entries (latest PIO update)
$data = @(
    @{ Row=24; A="MALCP-23"; B="4/14/2020"; C="Longos";       G=33;  H="F" },
    @{ Row=25; A="MALCP-24"; B="4/14/2020"; C="Hulong Duhat"; G=50;  H="F" },
    @{ Row=26; A="MALCP-25"; B="4/14/2020"; C="Tinajeros";    G=53;  H="F" },
    @{ Row=27; A="MALCP-26"; B="4/14/2020"; C="Tinajeros";    G=60;  H="M"; J=1; K="Dead" },
    @{ Row=28; A="MALCP-27"; B="4/14/2020"; C="Tugatog";      G="--"; H="M" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Formula = "=VLOOKUP(C$r,Barangay!`$A`$2:`$C`$22,3,0)"
    $ws.Cells.Item($r, 5).Formula = "=VLOOKUP(C$r,Barangay!`$A`$2:`$C`$22,2,0)"
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
    $ws.Cells.Item($r, 9).Value = 1
    if ($item.ContainsKey("J")) {
        $ws.Cells.Item($r, 10).Value = $item.J
    }
    if ($item.ContainsKey("K")) {
        $ws.Cells.Item($r, 11).Value = $item.K
    }
}

$ws.Range("A1").Worksheet.Application.Calculate()
